$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 503.7143
$ws.Range("I19").Value = 557.8333
$ws.Range("K19").Value = 557.8333
$ws.Range("M19").Value = -382.8333

$ws.Range("H33").Value = 711.6429000000001
$ws.Range("I33").Value = 733.36365
$ws.Range("J33").Value = 632
$ws.Range("K33").Value = 733.36365
$ws.Range("L33").Value = 632
$ws.Range("M33").Value = -504.36365
$ws.Range("N33").Value = -1090

$ws.Range("H86").Value = 1574.75
$ws.Range("I86").Value = 1199.5
$ws.Range("J86").Value = 1950
$ws.Range("K86").Value = 1199.5
$ws.Range("L86").Value = 1950
$ws.Range("M86").Value = -76.5
$ws.Range("N86").Value = -4196

$ws.Range("H89").Value = 1574.75
$ws.Range("I89").Value = 1199.5
$ws.Range("J89").Value = 1950
$ws.Range("K89").Value = 5997.5
$ws.Range("L89").Value = 9750
$ws.Range("M89").Value = -381.5
$ws.Range("N89").Value = -20982

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5770.909
$ws.Range("I32").Value = 5770.909
$ws.Range("K32").Value = 5770.909
$ws.Range("M32").Value = -5483.909

$ws.Range("H63").Value = 2000
$ws.Range("I63").Value = 2000
$ws.Range("K63").Value = 2000
$ws.Range("M63").Value = -1314

$ws.Range("H66").Value = 2000
$ws.Range("I66").Value = 2000
$ws.Range("K66").Value = 10000
$ws.Range("M66").Value = -6568

$ws.Range("H88").Value = 1721.7
$ws.Range("I88").Value = 1119.6666
$ws.Range("J88").Value = 2624.75
$ws.Range("K88").Value = 1119.6666
$ws.Range("L88").Value = 2624.75
$ws.Range("M88").Value = -713.6666
$ws.Range("N88").Value = -3436.75

$ws.Range("H91").Value = 1721.7
$ws.Range("I91").Value = 1119.6666
$ws.Range("J91").Value = 2624.75
$ws.Range("K91").Value = 1119.6666
$ws.Range("L91").Value = 2624.75
$ws.Range("M91").Value = 284.3334
$ws.Range("N91").Value = -5432.75

$ws.Range("H94").Value = 28931.4
$ws.Range("J94").Value = 28931.4
$ws.Range("L94").Value = 28931.4
$ws.Range("N94").Value = -30733.4

$ws.Range("H132").Value = 1161.0555
$ws.Range("I132").Value = 994.05884
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 2982.17652
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -452.17652
$ws.Range("N132").Value = -17060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 838.2
$ws.Range("I22").Value = 798
$ws.Range("K22").Value = 798
$ws.Range("M22").Value = -625

$ws.Range("H29").Value = 1220.0834
$ws.Range("I29").Value = 1480.3334
$ws.Range("J29").Value = 1133.3334
$ws.Range("K29").Value = 1480.3334
$ws.Range("L29").Value = 1133.3334
$ws.Range("M29").Value = -1191.3334
$ws.Range("N29").Value = -1711.3334

$ws.Range("H107").Value = 50297.5
$ws.Range("I107").Value = 57454.285
$ws.Range("J107").Value = 200
$ws.Range("K107").Value = 57454.285
$ws.Range("L107").Value = 200
$ws.Range("M107").Value = -55534.285
$ws.Range("N107").Value = -4040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 906.6667
$ws.Range("I22").Value = 898.2727
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 898.2727
$ws.Range("L22").Value = 999
$ws.Range("M22").Value = -548.2727
$ws.Range("N22").Value = -1699

$ws.Range("H47").Value = 36333.332

$ws.Range("H58").Value = 1586.375
$ws.Range("I58").Value = 1339.4
$ws.Range("K58").Value = 1339.4
$ws.Range("M58").Value = -1136.4

$ws.Range("H93").Value = 20371.25
$ws.Range("I93").Value = 8495
$ws.Range("K93").Value = 8495
$ws.Range("M93").Value = -6623

$ws.Range("H107").Value = 703.1429000000001
$ws.Range("I107").Value = 685.9091
$ws.Range("K107").Value = 685.9091
$ws.Range("M107").Value = 1234.0909

$ws.Range("H136").Value = 1586.375
$ws.Range("I136").Value = 1339.4
$ws.Range("K136").Value = 4018.2
$ws.Range("M136").Value = -1468.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 17.727272
$ws.Range("I7").Value = 17.222221
$ws.Range("K7").Value = 51.666663
$ws.Range("M7").Value = 60.333337

$ws.Range("H34").Value = 19990.818
$ws.Range("I34").Value = 149
$ws.Range("K34").Value = 447
$ws.Range("M34").Value = -363

$ws.Range("H108").Value = 362.6
$ws.Range("I108").Value = 362.6
$ws.Range("K108").Value = 1087.8
$ws.Range("M108").Value = 1792.2

$ws.Range("H114").Value = 1559.7142
$ws.Range("I114").Value = 1793.6
$ws.Range("K114").Value = 5380.799999999999
$ws.Range("M114").Value = -2126.799999999999

$ws.Range("H129").Value = 1105.4286
$ws.Range("J129").Value = 1733
$ws.Range("L129").Value = 5199
$ws.Range("N129").Value = -15199

$ws.Range("H131").Value = 4902
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 4902
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 14706
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -24786

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 11474.25
$ws.Range("I36").Value = 2949.5
$ws.Range("K36").Value = 2949.5
$ws.Range("M36").Value = -2464.5

$ws.Range("H43").Value = 3800
$ws.Range("I43").Value = 3800
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 3800
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -3649
$ws.Range("N43").ClearContents()

$ws.Range("H49").Value = 29919
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 29919
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 29919
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -30287

$ws.Range("H55").Value = 2000
$ws.Range("I55").Value = 2000
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 2000
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -1673
$ws.Range("N55").ClearContents()

$ws.Range("H95").Value = 6344
$ws.Range("J95").Value = 6344
$ws.Range("L95").Value = 6344
$ws.Range("N95").Value = -11836

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11806.75
$ws.Range("I7").Value = 11649.958
$ws.Range("J7").Value = 12747.5
$ws.Range("K7").Value = 11649.958
$ws.Range("L7").Value = 12747.5
$ws.Range("M7").Value = -11537.958
$ws.Range("N7").Value = -12971.5

$ws.Range("H55").Value = 549.5
$ws.Range("I55").Value = 411
$ws.Range("K55").Value = 411
$ws.Range("M55").Value = -238

$ws.Range("H69").Value = 75000
$ws.Range("J69").Value = 75000
$ws.Range("L69").Value = 75000
$ws.Range("N69").Value = -76622

$ws.Range("H72").Value = 75000
$ws.Range("J72").Value = 75000
$ws.Range("L72").Value = 225000
$ws.Range("N72").Value = -233112

$ws.Range("H82").Value = 3412.6365
$ws.Range("I82").Value = 3110
$ws.Range("J82").Value = 3664.8333
$ws.Range("K82").Value = 3110
$ws.Range("L82").Value = 3664.8333
$ws.Range("M82").Value = -2749
$ws.Range("N82").Value = -4386.8333

$ws.Range("H85").Value = 3412.6365
$ws.Range("I85").Value = 3110
$ws.Range("J85").Value = 3664.8333
$ws.Range("K85").Value = 3110
$ws.Range("L85").Value = 3664.8333
$ws.Range("M85").Value = -1862
$ws.Range("N85").Value = -6160.8333

$ws.Range("H93").Value = 41669132
$ws.Range("I93").Value = 47621356
$ws.Range("J93").Value = 3555
$ws.Range("K93").Value = 47621356
$ws.Range("L93").Value = 3555
$ws.Range("M93").Value = -47620108
$ws.Range("N93").Value = -6051

$ws.Range("H126").Value = 11806.75
$ws.Range("I126").Value = 11649.958
$ws.Range("J126").Value = 12747.5
$ws.Range("K126").Value = 34949.874
$ws.Range("L126").Value = 38242.5
$ws.Range("M126").Value = -32479.874
$ws.Range("N126").Value = -43182.5

$ws.Range("H127").Value = 100000
$ws.Range("J127").Value = 100000
$ws.Range("L127").Value = 100000
$ws.Range("N127").Value = -109920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 15333.667
$ws.Range("I62").Value = 6000
$ws.Range("J62").Value = 20000.5
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 20000.5
$ws.Range("M62").Value = -5376
$ws.Range("N62").Value = -21248.5

$ws.Range("H65").Value = 15333.667
$ws.Range("I65").Value = 6000
$ws.Range("J65").Value = 20000.5
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 100002.5
$ws.Range("M65").Value = -26880
$ws.Range("N65").Value = -106242.5

$ws.Range("H132").Value = 1314.8572
$ws.Range("I132").Value = 1121.5
$ws.Range("J132").Value = 2475
$ws.Range("K132").Value = 3364.5
$ws.Range("L132").Value = 7425
$ws.Range("M132").Value = -834.5
$ws.Range("N132").Value = -12485
